# Corrected dynamic case id in query for casedetails
# Remove the "caseDetailQuery" column (C) entirely - both its header and
# its query text (which contained a hard-coded case id). Columns D and E
# shift left to become C and D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire column C (caseDetailQuery header + its query text).
# This shifts columns D:E left to C:D automatically.
$ws.Range("C:C").Delete() | Out-Null

# Update the active selection to match the post-edit state (B2).
$ws.Range("B2").Select() | Out-Null
